$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 becomes the "Gold Feb 26 / GC=F" futures row
$ws.Range("B2").Value = "Gold Feb 26"
$ws.Range("C2").Value = "GC=F"
$ws.Range("D2").Value = 4238.5
$ws.Range("E2").Value = 55.7
$ws.Range("F2").Value = 1.76
$ws.Range("G2").Value = 60
$ws.Range("H2").Value = 63
$ws.Range("I2").Value = 83
$ws.Range("J2").Value = 83
$ws.Range("K2").Value = 67.59999999999999
$ws.Range("N2").Value = 54.82400714602223

# Row 3 becomes the "StreetTRACKS Gold Shares / GLD" row
$ws.Range("B3").Value = "StreetTRACKS Gold Shares"
$ws.Range("C3").Value = "GLD"
$ws.Range("D3").Value = 387.13
$ws.Range("E3").Value = 56.3
$ws.Range("F3").Value = 1.05
$ws.Range("G3").Value = 60
$ws.Range("H3").Value = 73
$ws.Range("I3").Value = 83
$ws.Range("J3").Value = 96
$ws.Range("K3").Value = 67.59999999999999
$ws.Range("N3").Value = 54.82400714602223

# Row 4 (Newmont / NEM) keeps same B/C/D-J, only K and N update
$ws.Range("K4").Value = 66.40000000000001
$ws.Range("N4").Value = 54.82400714602223
